$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the mimetype cells for the two tab-delimited text log files (F14, F16)
# from "text/x.vnd.abi.plot+csv" to the new, more specific mimetype.
$ws.Range("F14").Value = "text/x.vnd.abi.plot+Tab-separated-values"
$ws.Range("F16").Value = "text/x.vnd.abi.plot+Tab-separated-values"

# Update the active selection on the sheet to G28 (per saved view state).
$ws.Range("G28").Select()
